$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the summary values at the top of the statement
$ws.Range("E11").Value = 53282
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Roll the last worker's data row up into the first (and only remaining)
# detail row before deleting the now-redundant repeated rows below it.
$ws.Range("C16").Value = "1103114066"
$ws.Range("D16").Value = "ANDREA CAROLINA PEREZ CHICA"
$ws.Range("E16").Value = "2204"
$ws.Range("F16").Value = 53282
$ws.Range("G16").Value = 1332045

# Remove the now-duplicated detail rows (17-21 repeated the same worker,
# row 22 held the data that was just copied up into row 16). Deleting this
# block shifts the signature rows (old 27/28) up to 21/22.
$ws.Rows("17:22").Delete()
